# "finish dev of barracks"
# Insert a new "soldierMaterials" worksheet (soldier equipment costs table)
# right before the "houses" sheet, modelled on the existing "materials"
# sheet (same A1:I2 layout / styling), then populate it with the new
# soldier-equipment item columns.

$wb = $excel.ActiveWorkbook

$materialsSheet = $wb.Worksheets.Item("materials")
$housesSheet    = $wb.Worksheets.Item("houses")

# Update the houses sheet's view state: it will no longer be the
# tab-selected sheet, and its cursor moves to D5.
$housesSheet.Activate()
$housesSheet.Range("D5").Select()

# Create the new sheet directly before "houses" (reuses the rId/sheet
# file that "houses" used to occupy; houses gets a freshly allocated one).
$newSheet = $wb.Worksheets.Add($housesSheet)
$newSheet.Name = "soldierMaterials"

# The new sheet inherits its own copies of the legacy (REF!-broken)
# AutoFilter built-in names, same as "houses" carries.
$newSheet.Names.Add("Excel_BuiltIn__FilterDatabase_2", "=#REF!")
$newSheet.Names.Add("Excel_BuiltIn__FilterDatabase_6", "=#REF!")
$newSheet.Names.Add("Excel_BuiltIn__FilterDatabase_7", "=#REF!")
$newSheet.Names.Add("Excel_BuiltIn__FilterDatabase_8", "=#REF!")

# Clone layout/formatting from "materials" (same header/data row styles).
$materialsSheet.Range("A1:I2").Copy($newSheet.Range("A1:I2"))

# Overwrite the header labels with the soldier-equipment fields.
$newSheet.Range("A1").Value = "INT_level"
$newSheet.Range("B1").Value = "INT_deathHand"
$newSheet.Range("C1").Value = "INT_heroBones"
$newSheet.Range("D1").Value = "INT_soulStone"
$newSheet.Range("E1").Value = "INT_magicBox"
$newSheet.Range("F1").Value = "INT_confessionHood"
$newSheet.Range("G1").Value = "INT_brightRing"
$newSheet.Range("H1").Value = "INT_holyBook"
$newSheet.Range("I1").Value = "INT_brightAlloy"

# Data row.
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2:I2").Value = 1000

$newSheet.Activate()
$newSheet.Range("H3").Select()
